# Automatische test-sync: 2025-08-26 21:07:50
# Append a new log row (row 8) to the "Logs" sheet and update the
# "Retour / Terugbetaling" tally on the "Dashboard" sheet.

$wb = $excel.ActiveWorkbook

# --- Logs sheet -----------------------------------------------------------
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A8").Value = "Retour status"
$logs.Range("B8").Value = "mailmind.test@zohomail.eu"
$logs.Range("D8").Value = "Retour / Terugbetaling"
$logs.Range("F8").Value = "2025-08-26 21:07:17"
$logs.Range("G8").Value = "Nee"
$logs.Range("H8").Value = "Ja"
$logs.Range("I8").Value = "Nee"
$logs.Range("J8").Value = "Nee"

# Extend the conditional formatting ranges so the new row is covered too.
$logs.Range("D2:D7").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("D2:D8"))
$logs.Range("G2:G7").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("G2:G8"))
$logs.Range("H2:H7").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("H2:H8"))
$logs.Range("I2:I7").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("I2:I8"))
$logs.Range("J2:J7").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("J2:J8"))

# --- Dashboard sheet --------------------------------------------------------
$dashboard = $wb.Worksheets.Item("Dashboard")

# "Retour / Terugbetaling" count goes up from 2 to 3.
$dashboard.Range("B3").Value = 3
